$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the "Meta description: ..." paragraph that currently sits
#    right after the title heading paragraph.
# ------------------------------------------------------------------
$pCount = $d.Paragraphs.Count
for ($i = 1; $i -le $pCount; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Meta description*") {
        [void]$p.Range.Delete()
        break
    }
}

# ------------------------------------------------------------------
# 2) At the end of the document, split the final (italic) paragraph so a
#    new bold paragraph ("Play Ultra Hot Deluxe Free - ...") precedes
#    it, and update the italic paragraph's text to the new copy.
# ------------------------------------------------------------------
$pCount = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($pCount)
$r = $lastPara.Range
$rng = $d.Range($r.Start, $r.End - 1)

$snippet = @"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Ultra Hot Deluxe Free - Simple Gameplay and Additional Bet Game</w:t></w:r></w:p>
<w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of Ultra Hot Deluxe, a basic online slot game with an additional 50-50 bet game. Play for free on any device.</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

[void]$rng.InsertXML($snippet)

Write-Host "Done. ParaCount:" $d.Paragraphs.Count
